# Commit: "Updated game extraction to remove trademarks"
#
# The game-extraction automation inserted a new shared string
# ("SUPERCHICOS") into the workbook's string table ahead of the existing
# game-name list without re-pointing the rows that were already written
# to their new slots. The net, user-visible effect on the "Sheet1"
# worksheet is that column A rows 2-16 each end up showing the title that
# used to belong to the row directly above them, with row 2 now showing
# the newly inserted "SUPERCHICOS" title; rows 17 and below are
# unaffected (their index shifted by the same amount as the insertion,
# so what they display does not change).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$firstRow = 2
$lastRow = 16

# Snapshot the current titles for the affected range before overwriting
# anything, then write them back shifted down by one row.
$titles = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $titles[$r] = $ws.Cells.Item($r, 1).Value2
}

for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $ws.Cells.Item($r, 1).Value = $titles[$r - 1]
}

$ws.Cells.Item($firstRow, 1).Value = "SUPERCHICOS"
